# Auto-generated edit script applying the Betfair odds update for 2025-12-18
# Sets updated numeric values for the changed cells (rows 2-19) per the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("AA2").Value = 180
$ws.Range("AD2").Value = 24
$ws.Range("AI2").Value = 70
$ws.Range("AL2").Value = 30
$ws.Range("AN2").Value = 7.2
$ws.Range("AO2").Value = 75
$ws.Range("G2").Value = 1.67
$ws.Range("H2").Value = 5.7
$ws.Range("J2").Value = 4.2
$ws.Range("P2").Value = 2.34
$ws.Range("T2").Value = 1.72
$ws.Range("U2").Value = 2.16
$ws.Range("X2").Value = 23
# Row 3
$ws.Range("I3").Value = 3.4
$ws.Range("J3").Value = 3.55
$ws.Range("P3").Value = 2.12
$ws.Range("Q3").Value = 1.75
# Row 4
$ws.Range("AB4").Value = 11.5
$ws.Range("AC4").Value = 22
$ws.Range("AD4").Value = 1000
$ws.Range("AE4").Value = 320
$ws.Range("AF4").Value = 9.4
$ws.Range("AG4").Value = 13.5
$ws.Range("AH4").Value = 1000
$ws.Range("AJ4").Value = 9.6
$ws.Range("AK4").Value = 15
$ws.Range("AL4").Value = 55
$ws.Range("AM4").Value = 240
$ws.Range("AN4").Value = 4.2
$ws.Range("F4").Value = 1.24
$ws.Range("G4").Value = 1.27
$ws.Range("H4").Value = 13.5
$ws.Range("I4").Value = 17.5
$ws.Range("J4").Value = 7.4
$ws.Range("K4").Value = 8.199999999999999
$ws.Range("M4").Value = 1.02
$ws.Range("N4").Value = 5.8
$ws.Range("Q4").Value = 1.55
$ws.Range("R4").Value = 1.7
$ws.Range("T4").Value = 2.2
$ws.Range("X4").Value = 38
$ws.Range("Z4").Value = 180
# Row 5
$ws.Range("G5").Value = 1.71
$ws.Range("J5").Value = 3.85
$ws.Range("K5").Value = 4.6
$ws.Range("M5").Value = 1.08
$ws.Range("P5").Value = 1.74
$ws.Range("Q5").Value = 2.08
# Row 6
$ws.Range("AN6").Value = 17
$ws.Range("I6").Value = 4.7
$ws.Range("N6").Value = 3.45
$ws.Range("X6").Value = 17
$ws.Range("Y6").Value = 980
# Row 7
$ws.Range("G7").Value = 1.34
$ws.Range("H7").Value = 9.4
$ws.Range("I7").Value = 15
$ws.Range("J7").Value = 5.5
$ws.Range("K7").Value = 6.8
$ws.Range("P7").Value = 2.28
$ws.Range("Q7").Value = 1.59
$ws.Range("T7").Value = 2.16
$ws.Range("U7").Value = 1.74
$ws.Range("Y7").Value = 55
# Row 8
$ws.Range("AA8").Value = 130
$ws.Range("P8").Value = 1.7
# Row 9
$ws.Range("F9").Value = 2.56
$ws.Range("H9").Value = 2.8
$ws.Range("P9").Value = 1.8
$ws.Range("Q9").Value = 1.88
# Row 10
$ws.Range("G10").Value = 3.55
$ws.Range("H10").Value = 2.72
$ws.Range("I10").Value = 2.8
$ws.Range("J10").Value = 2.82
# Row 11
$ws.Range("AH11").Value = 1000
$ws.Range("F11").Value = 1.22
$ws.Range("G11").Value = 1.27
$ws.Range("H11").Value = 12.5
$ws.Range("I11").Value = 15.5
$ws.Range("J11").Value = 7.2
$ws.Range("K11").Value = 8.800000000000001
$ws.Range("P11").Value = 3.1
$ws.Range("T11").Value = 1.9
$ws.Range("U11").Value = 1.9
# Row 12
$ws.Range("AB12").Value = 1000
$ws.Range("AE12").Value = 1000
$ws.Range("AF12").Value = 8.6
$ws.Range("AG12").Value = 14.5
$ws.Range("AI12").Value = 210
$ws.Range("AM12").Value = 240
$ws.Range("AN12").Value = 5
$ws.Range("F12").Value = 1.28
$ws.Range("H12").Value = 11
$ws.Range("I12").Value = 14
$ws.Range("J12").Value = 6.2
$ws.Range("P12").Value = 2.42
$ws.Range("Q12").Value = 1.56
$ws.Range("T12").Value = 2.26
$ws.Range("U12").Value = 1.6
$ws.Range("Z12").Value = 160
# Row 13
$ws.Range("AA13").Value = 480
$ws.Range("AB13").Value = 8.199999999999999
$ws.Range("AC13").Value = 11
$ws.Range("AD13").Value = 1000
$ws.Range("AE13").Value = 180
$ws.Range("AG13").Value = 11
$ws.Range("AI13").Value = 160
$ws.Range("AK13").Value = 18
$ws.Range("AL13").Value = 48
$ws.Range("AM13").Value = 210
$ws.Range("F13").Value = 1.46
$ws.Range("G13").Value = 1.51
$ws.Range("H13").Value = 7.8
$ws.Range("J13").Value = 4.3
$ws.Range("P13").Value = 1.91
$ws.Range("Q13").Value = 1.8
$ws.Range("T13").Value = 2.12
$ws.Range("U13").Value = 1.72
$ws.Range("Y13").Value = 29
$ws.Range("Z13").Value = 95
# Row 14
$ws.Range("F14").Value = 1.19
$ws.Range("H14").Value = 11.5
$ws.Range("I14").Value = 32
$ws.Range("P14").Value = 2.6
$ws.Range("Q14").Value = 1.5
$ws.Range("Z14").Value = 350
# Row 15
$ws.Range("AH15").Value = 18.5
$ws.Range("AI15").Value = 55
$ws.Range("I15").Value = 3.75
$ws.Range("P15").Value = 1.92
$ws.Range("Q15").Value = 1.87
$ws.Range("T15").Value = 1.71
$ws.Range("U15").Value = 2.12
# Row 16
$ws.Range("F16").Value = 1.3
$ws.Range("G16").Value = 1.31
$ws.Range("K16").Value = 8
$ws.Range("P16").Value = 2.9
$ws.Range("T16").Value = 1.91
$ws.Range("U16").Value = 1.9
$ws.Range("X16").Value = 1000
# Row 17
$ws.Range("AC17").Value = 8.800000000000001
$ws.Range("F17").Value = 3.45
$ws.Range("P17").Value = 2.04
$ws.Range("Q17").Value = 1.8
$ws.Range("T17").Value = 1.66
# Row 18
$ws.Range("F18").Value = 1.88
$ws.Range("G18").Value = 2
$ws.Range("H18").Value = 4.3
$ws.Range("I18").Value = 5.1
$ws.Range("J18").Value = 3.4
$ws.Range("K18").Value = 3.95
$ws.Range("M18").Value = 1.08
$ws.Range("P18").Value = 1.77
$ws.Range("Q18").Value = 1.92
$ws.Range("T18").Value = 1.88
$ws.Range("U18").Value = 1.94
# Row 19
$ws.Range("AA19").Value = 160
$ws.Range("F19").Value = 1.67
$ws.Range("G19").Value = 1.81
$ws.Range("H19").Value = 5
$ws.Range("I19").Value = 5.8
$ws.Range("K19").Value = 4.5
$ws.Range("P19").Value = 2.1
$ws.Range("Q19").Value = 1.74
$ws.Range("U19").Value = 2.1
$ws.Range("X19").Value = 1000
